# TrialsSetup.xlsx update (2026-01-19 16:00)
#
# The underlying Power Query ("Query1", refreshed from the SharePoint-backed
# "Set-up" sheet) dropped the "PKN605" trial row entirely and the "QUILT"
# trial's Progress value moved from 25% to 37.5%. Reproduce that refreshed
# state directly against the worksheet/table/defined name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "PKN605" was row 13 (Trial Name | Days remaining | Progress). Removing the
# row shifts ALPINE / REMASTER (CLOU) up one row, shrinks the table to
# A1:C14, and re-numbers the shared-string table automatically on save.
$ws.Rows(13).Delete()

# QUILT's Progress (row 3, column C) was refreshed from 25 to 37.5.
$ws.Range("C3").Value = 37.5

# Keep the workbook-level "ExternalData_1" defined name (used by the query
# table) in sync with the new, smaller data extent.
$wb.Names("ExternalData_1").RefersTo = "=Sheet1!`$A`$1:`$C`$14"
